# Daily attendance processing - 2026-01-24 13:47:12
#
# Normalizes the "Recorded By" column (G) on the active sheet so that the
# comma-separated list of recorders in each cell follows a consistent
# ordering/priority instead of the arbitrary order they were originally
# appended in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Priority order used to sort the comma separated "Recorded By" names.
# Lower value = sorts earlier. Anything not found here keeps its relative
# order and is placed after all known names (stable sort). Comparison must
# be case-sensitive (e.g. "system" and "System" are different recorders),
# so a plain PowerShell hashtable (case-insensitive keys) cannot be used.
function Get-RecorderPriority($name) {
    if ($name.Equals("system")) { return 0 }
    if ($name.Equals("dnasr281@gmail.com")) { return 1 }
    if ($name.Equals("admin@admin.com")) { return 2 }
    if ($name.Equals("backup@backdoor.com")) { return 3 }
    if ($name.Equals("System")) { return 4 }
    return 999
}

# Determine the last used row in column G (Recorded By).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $items = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $name = $parts[$i]
        $pr = Get-RecorderPriority $name
        $items += [PSCustomObject]@{ Pri = $pr; Idx = $i; Name = $name }
    }

    $sorted = $items | Sort-Object Pri, Idx
    $newParts = $sorted | ForEach-Object { $_.Name }
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
